# Auto-generated edit script: refresh market price / profit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2160.1667
$ws.Range("J17").Value = 2160.1667
$ws.Range("L17").Value = 6480.500100000001
$ws.Range("N17").Value = -6816.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1051.375
$ws.Range("I2").Value = 1002.75
$ws.Range("K2").Value = 1002.75
$ws.Range("M2").Value = -889.75
$ws.Range("H32").Value = 17248.633
$ws.Range("I32").Value = 13154.417
$ws.Range("J32").Value = 33625.5
$ws.Range("K32").Value = 13154.417
$ws.Range("L32").Value = 33625.5
$ws.Range("M32").Value = -12867.417
$ws.Range("N32").Value = -34199.5
$ws.Range("H61").Value = 3156.743
$ws.Range("I61").Value = 1756.619
$ws.Range("J61").Value = 5256.9287
$ws.Range("K61").Value = 1756.619
$ws.Range("L61").Value = 5256.9287
$ws.Range("M61").Value = -1544.619
$ws.Range("N61").Value = -5680.9287
$ws.Range("H74").Value = 2761.4075
$ws.Range("I74").Value = 2503.2
$ws.Range("K74").Value = 2503.2
$ws.Range("M74").Value = -1629.2
$ws.Range("H77").Value = 2761.4075
$ws.Range("I77").Value = 2503.2
$ws.Range("K77").Value = 12516
$ws.Range("M77").Value = -8148
$ws.Range("H88").Value = 10443842
$ws.Range("I88").Value = 23139.6
$ws.Range("J88").Value = 27811680
$ws.Range("K88").Value = 23139.6
$ws.Range("L88").Value = 27811680
$ws.Range("M88").Value = -22733.6
$ws.Range("N88").Value = -27812492
$ws.Range("H91").Value = 10443842
$ws.Range("I91").Value = 23139.6
$ws.Range("J91").Value = 27811680
$ws.Range("K91").Value = 23139.6
$ws.Range("L91").Value = 27811680
$ws.Range("M91").Value = -21735.6
$ws.Range("N91").Value = -27814488
$ws.Range("H103").Value = 70000
$ws.Range("J103").Value = 70000
$ws.Range("L103").Value = 70000
$ws.Range("N103").Value = -72344
$ws.Range("H110").Value = 3432.625
$ws.Range("I110").Value = 3340.1365
$ws.Range("K110").Value = 3340.1365
$ws.Range("M110").Value = -1295.1365
$ws.Range("H116").Value = 1051.375
$ws.Range("I116").Value = 1002.75
$ws.Range("K116").Value = 1002.75
$ws.Range("M116").Value = 1291.25
$ws.Range("H136").Value = 3156.743
$ws.Range("I136").Value = 1756.619
$ws.Range("J136").Value = 5256.9287
$ws.Range("K136").Value = 5269.857
$ws.Range("L136").Value = 15770.7861
$ws.Range("M136").Value = -2719.857
$ws.Range("N136").Value = -20870.7861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1051.375
$ws.Range("I3").Value = 1002.75
$ws.Range("K3").Value = 1002.75
$ws.Range("M3").Value = -888.75
$ws.Range("H105").Value = 3927.4688
$ws.Range("I105").Value = 3030.8696
$ws.Range("K105").Value = 3030.8696
$ws.Range("M105").Value = -1283.8696
$ws.Range("H111").Value = 32666.666
$ws.Range("J111").Value = 32666.666
$ws.Range("L111").Value = 32666.666
$ws.Range("N111").Value = -40846.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4172.371
$ws.Range("I31").Value = 3744.6875
$ws.Range("J31").Value = 4321.1304
$ws.Range("K31").Value = 3744.6875
$ws.Range("L31").Value = 4321.1304
$ws.Range("M31").Value = -3449.6875
$ws.Range("N31").Value = -4911.1304
$ws.Range("H34").Value = 4172.371
$ws.Range("I34").Value = 3744.6875
$ws.Range("J34").Value = 4321.1304
$ws.Range("K34").Value = 3744.6875
$ws.Range("L34").Value = 4321.1304
$ws.Range("M34").Value = -3542.6875
$ws.Range("N34").Value = -4725.1304
$ws.Range("H132").Value = 5544.75
$ws.Range("I132").Value = 5958.5713
$ws.Range("J132").Value = 4965.4
$ws.Range("K132").Value = 17875.7139
$ws.Range("L132").Value = 14896.2
$ws.Range("M132").Value = -15345.7139
$ws.Range("N132").Value = -19956.2
$ws.Range("H134").Value = 4550.8076
$ws.Range("I134").Value = 2856
$ws.Range("J134").Value = 7752.1113
$ws.Range("K134").Value = 8568
$ws.Range("L134").Value = 23256.3339
$ws.Range("M134").Value = -6033
$ws.Range("N134").Value = -28326.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 202224.5
$ws.Range("J57").Value = 202224.5
$ws.Range("L57").Value = 606673.5
$ws.Range("N57").Value = -607791.5
$ws.Range("H122").Value = 1418.6364
$ws.Range("I122").Value = 1112.2858
$ws.Range("J122").Value = 1561.6
$ws.Range("K122").Value = 10010.5722
$ws.Range("L122").Value = 14054.4
$ws.Range("M122").Value = -7560.572200000001
$ws.Range("N122").Value = -18954.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1983.3864
$ws.Range("I102").Value = 1506.5952
$ws.Range("K102").Value = 1506.5952
$ws.Range("M102").Value = 115.4048
$ws.Range("H122").Value = 4734.381
$ws.Range("I122").Value = 1799.3636
$ws.Range("J122").Value = 7962.9
$ws.Range("K122").Value = 5398.0908
$ws.Range("L122").Value = 23888.7
$ws.Range("M122").Value = -2948.0908
$ws.Range("N122").Value = -28788.7
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 3989.6
$ws.Range("I126").Value = 1316
$ws.Range("K126").Value = 3948
$ws.Range("M126").Value = -1478
$ws.Range("H128").Value = 60000
$ws.Range("J128").Value = 60000
$ws.Range("L128").Value = 60000
$ws.Range("N128").Value = -69960
$ws.Range("H132").Value = 3649.66
$ws.Range("I132").Value = 3292.2683
$ws.Range("K132").Value = 9876.804900000001
$ws.Range("M132").Value = -7346.804900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 56666.332
$ws.Range("J101").Value = 56666.332
$ws.Range("L101").Value = 56666.332
$ws.Range("N101").Value = -63156.332
$ws.Range("H106").Value = 20500
$ws.Range("J106").Value = 20500
$ws.Range("L106").Value = 20500
$ws.Range("N106").Value = -23024
$ws.Range("H136").Value = 3723.652
$ws.Range("I136").Value = 2689.0303
$ws.Range("K136").Value = 8067.090899999999
$ws.Range("M136").Value = -5517.090899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 24186.137
$ws.Range("I96").Value = 39762.31
$ws.Range("J96").Value = 1687.2222
$ws.Range("K96").Value = 39762.31
$ws.Range("L96").Value = 1687.2222
$ws.Range("M96").Value = -38389.31
$ws.Range("N96").Value = -4433.2222
$ws.Range("H107").Value = 633.44446
$ws.Range("I107").Value = 557.2857
$ws.Range("K107").Value = 1671.8571
$ws.Range("M107").Value = 248.1428999999998
$ws.Range("H126").Value = 2775.7
$ws.Range("I126").Value = 2692.4285
$ws.Range("K126").Value = 8077.2855
$ws.Range("M126").Value = -5607.2855
